$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill column C (Sub-Category) with "iGaming" for all data rows (2 through 437),
# mirroring the value already present in column B (Category).
$ws.Range("C2:C437").Value = "iGaming"

# Update the active selection to match the final cursor position recorded in the diff.
$ws.Range("C8").Select()
